# Replace the half-width square brackets around the English dialogue
# text with full-width brackets 【 】 for the four updated lines, keeping
# everything else (including the embedded newline) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[PopupDialog(dialogHead=`"`$avatar_amiya2`")] 【Kashchey... Your plot must be stopped. Right here.】`n"
$ws.Range("C7").Value = "[PopupDialog(dialogHead=`"`$avatar_chen2`")] 【Clear Eyes】`n"
$ws.Range("C8").Value = "[PopupDialog(dialogHead=`"`$avatar_amiya2`")] 【It’s not this sword that will shatter your plans, Kashchey.】`n"
$ws.Range("C9").Value = "[PopupDialog(dialogHead=`"`$avatar_amiya2`")] 【It’s the person holding the sword.】`n"
